$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 11:52"

$ws.Range("B4").Value = 61171
$ws.Range("C4").Value = 36707
$ws.Range("D4").Value = 16288
$ws.Range("E4").Value = 8176

$ws.Range("B5").Value = 48916
$ws.Range("C5").Value = 19615
$ws.Range("D5").Value = 24326
$ws.Range("E5").Value = 4975

$ws.Range("B6").Value = 16885
$ws.Range("C6").Value = 6569
$ws.Range("D6").Value = 8564
$ws.Range("E6").Value = 1752

$ws.Range("B7").Value = 15832
$ws.Range("C7").Value = 5512
$ws.Range("D7").Value = 7857
$ws.Range("E7").Value = 2463

$ws.Range("B8").Value = 12701
$ws.Range("C8").Value = 10936
$ws.Range("D8").Value = 469
$ws.Range("E8").Value = 1296

$ws.Range("B9").Value = 12048
$ws.Range("C9").Value = 5710
$ws.Range("D9").Value = 5131
$ws.Range("E9").Value = 1207

$ws.Range("B10").Value = 8697
$ws.Range("C10").Value = 5573
$ws.Range("D10").Value = 2577

$ws.Range("B14").Value = 5091
$ws.Range("C14").Value = 2291
$ws.Range("D14").Value = 2061
$ws.Range("E14").Value = 739

$ws.Range("B15").Value = 4815
$ws.Range("C15").Value = 2113
$ws.Range("D15").Value = 2251
$ws.Range("E15").Value = 451

$ws.Range("B17").Value = 3918
$ws.Range("C17").Value = 2159
$ws.Range("D17").Value = 1427
$ws.Range("E17").Value = 332

$ws.Range("B23").Value = 2785
$ws.Range("C23").Value = 1887
$ws.Range("D23").Value = 452
$ws.Range("E23").Value = 446

$ws.Range("B30").Value = 2283
$ws.Range("C30").Value = 835
$ws.Range("D30").Value = 1175
$ws.Range("E30").Value = 273

$ws.Range("B32").Value = 2205
$ws.Range("C32").Value = 1131
$ws.Range("D32").Value = 939
$ws.Range("E32").Value = 135

$ws.Range("B33").Value = 2173
$ws.Range("C33").Value = 1483
$ws.Range("D33").Value = 499

$ws.Range("B38").Value = 1486
$ws.Range("C38").Value = 1202
$ws.Range("D38").Value = 154

$ws.Range("C59").Value = 96
$ws.Range("D59").Value = 16
